$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
$ws.Range("D2").Value = "'43.940.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.43%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.281.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.39%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'234.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.14%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +3.88%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'65.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +9.14%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.432"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +7.37%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +17.63%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'57.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.97%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'26.39"
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.70%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.619.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.31%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'15.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.50%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +5.45%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.833"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.89%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.282.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.30%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'43.695.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.70%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0997"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +10.81%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'74.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.79%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +2.02%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'263.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +8.55%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E25").Value = "'  +6.58%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.03%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'10.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.64%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'172.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.01%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'21.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +7.01%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -1.86%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.65%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +8.11%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.00%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +6.14%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.79%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'4.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.90%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'RenderToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'3.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +9.25%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'THORChain"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'6.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +8.44%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  +0.33%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +4.21%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.19%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'8.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.65%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'17.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.11%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.0981"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.83%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'FTXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'4.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.74%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'98.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.55%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Celestia"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'10.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +22.13%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.70%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.479.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.03%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'2.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +7.18%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.000207"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -12.96%  "
$ws.Range("E51").Style = "Normal"
